$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit performs a cyclic rotation of the "observation" data among rows 7, 8
# and 9 (columns A, B, E, F, G, H, I, Q, R, AC):
#   new row7 = old row9
#   new row8 = old row7
#   new row9 = old row8
# All other columns (C, D, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW,
# AX, AY) are identical across the three rows already, so nothing to do there.

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AC")

# Snapshot the current ("before") values for the rows involved, since the
# rotation needs the original data of all three rows simultaneously.
$row7 = @{}
$row8 = @{}
$row9 = @{}
foreach ($c in $cols) {
    $row7[$c] = $ws.Range($c + "7").Value2
    $row8[$c] = $ws.Range($c + "8").Value2
    $row9[$c] = $ws.Range($c + "9").Value2
}
# Column I holds numeric-looking text (e.g. "20", "15", or blank) that must
# stay stored as text, so grab it via .Text to keep it a string.
$row7["I"] = $ws.Range("I7").Text
$row8["I"] = $ws.Range("I8").Text
$row9["I"] = $ws.Range("I9").Text

function Set-RowValues($targetRow, $data) {
    foreach ($c in $cols) {
        $ws.Range($c + $targetRow).Value2 = $data[$c]
    }
    # Write column I explicitly as text, preserving the default (unstyled)
    # cell format -- matches how these values are stored in the source file.
    $cell = $ws.Range("I" + $targetRow)
    if ([string]::IsNullOrEmpty($data["I"])) {
        $cell.ClearContents()
    } else {
        $cell.NumberFormat = "@"
        $cell.Value = $data["I"]
        $cell.Style = "Normal"
    }
}

Set-RowValues 7 $row9
Set-RowValues 8 $row7
Set-RowValues 9 $row8
